# Update the "FSS" worksheet to document / demonstrate multi-FSS and
# multi-MountTarget scenarios (commit: "updated FSS code to inlude multi
# FSS and multi MT scenarios").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FSS")

# --- Row 1: updated instructions banner, taller to fit the extra text ---
$ws.Range("A1").Value = "Columns: Region, Compartment Name, Availability Domain, MountTarget Name, MountTarget SubnetName, FSS name, Path are mandatory.`nDefault value for columns if left blank: sourceCIDR- 0.0.0.0/0, Access- READ_ONLY, GID- 65534, UID- 65534, IDSquash- NONE and require_ps_port- false`nMount target IP, FSS Capacity, FSS Inodes will take default values from OCI if left blank`nResources will be created based on MountTargetName and FSSName columns.`nBelow sample data shows example of  multiple FSS(FSS1 and FSS2) using single MT(MT1) , 1 FSS(FSS3) using multiple MTs(MT2 and MT3) and also 1 FSS(FSS4) using 1 MT(MT4)"
$ws.Rows.Item(1).RowHeight = 84.5

# --- Prepare a fresh row 8 to host the old "<END>" marker row, formatted
# like the rest of the data rows (thin-bordered cells) ---
$ws.Cells.Item(8, 1).Value = "<END>"
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)

# --- Row 3: MT1 / FSS1 in AD2, subnet HUB_SERVICES_SN, path /fss1/ ---
$ws.Cells.Item(3, 2).Value = "OCICtoOCI"
$ws.Cells.Item(3, 3).Value = "AD2"
$ws.Cells.Item(3, 4).Value = "MT1"
$ws.Cells.Item(3, 5).Value = "HUB_SERVICES_SN"
$ws.Cells.Item(3, 10).Value = "FSS1"
$ws.Cells.Item(3, 11).Value = "/fss1/"

# --- Row 4: MT1 / FSS2 in AD2, subnet HUB_SERVICES_SN, path /fss2/ ---
# (second FSS sharing the same Mount Target as row 3)
$ws.Cells.Item(4, 1).Value = "Ashburn"
$ws.Cells.Item(4, 2).Value = "OCICtoOCI"
$ws.Cells.Item(4, 3).Value = "AD2"
$ws.Cells.Item(4, 4).Value = "MT1"
$ws.Cells.Item(4, 5).Value = "HUB_SERVICES_SN"
$ws.Cells.Item(4, 10).Value = "FSS2"
$ws.Cells.Item(4, 11).Value = "/fss2/"

# --- Row 5: MT3 / FSS3 in AD3, subnet HUB_SERVICES_SN, path /fss3/,
# export options sourceCIDR 11.0.0.0/8, Access READ_WRITE ---
$ws.Cells.Item(5, 1).Value = "Ashburn"
$ws.Cells.Item(5, 2).Value = "OCICtoOCI"
$ws.Cells.Item(5, 3).Value = "AD3"
$ws.Cells.Item(5, 4).Value = "MT3"
$ws.Cells.Item(5, 5).Value = "HUB_SERVICES_SN"
$ws.Cells.Item(5, 10).Value = "FSS3"
$ws.Cells.Item(5, 11).Value = "/fss3/"
$ws.Cells.Item(5, 12).Value = "11.0.0.0/8"
$ws.Cells.Item(5, 13).Value = "READ_WRITE"

# --- Row 6: MT2 / FSS3 in AD3 (same FSS reachable from a second Mount
# Target), export options sourceCIDR 10.0.0.0/8, Access READ_ONLY ---
$ws.Cells.Item(6, 1).Value = "Ashburn"
$ws.Cells.Item(6, 2).Value = "OCICtoOCI"
$ws.Cells.Item(6, 3).Value = "AD3"
$ws.Cells.Item(6, 4).Value = "MT2"
$ws.Cells.Item(6, 5).Value = "HUB_SERVICES_SN"
$ws.Cells.Item(6, 10).Value = "FSS3"
$ws.Cells.Item(6, 11).Value = "/fss3/"
$ws.Cells.Item(6, 12).Value = "10.0.0.0/8"
$ws.Cells.Item(6, 13).Value = "READ_ONLY"

# --- Row 7: MT4 / FSS4 in AD3, subnet HUB_SERVICES_SN, path /fss4/ ---
$ws.Cells.Item(7, 1).Value = "Ashburn"
$ws.Cells.Item(7, 2).Value = "OCICtoOCI"
$ws.Cells.Item(7, 3).Value = "AD3"
$ws.Cells.Item(7, 4).Value = "MT4"
$ws.Cells.Item(7, 5).Value = "HUB_SERVICES_SN"
$ws.Cells.Item(7, 10).Value = "FSS4"
$ws.Cells.Item(7, 11).Value = "/fss4/"
